# Generate Report for Handoff
# The "b.md" row now has a fresh handoff in flight for both locales, so its
# status flips from "Handed back: in sync with en-US" to "Ready for handoff"
# and the Latest Handoff File / Latest Handoff Datetime columns point at the
# newly generated xlf + timestamp.

$wb = $excel.ActiveWorkbook

function Set-HandoffCell($ws, $cellRef, $newText) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $newText
    $rngAddr = $rng.Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $rngAddr) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet: b.md row (row 3) status -> Ready for handoff ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
Set-HandoffCell $wsZh "C3" "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-10 06:44:48"

# --- de-de sheet: b.md row (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
Set-HandoffCell $wsDe "C3" "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-10 06:44:54"
